$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the current column F ("Comments"), shifting
# the old F column (and everything after it) two places to the right.
$ws.Columns("F:G").Insert()

# New column headers for the inserted columns.
$ws.Range("F1").Value = "Variable"
$ws.Range("G1").Value = "Variable type"

# Fill in the "Variable" / "Variable type" columns row by row (top to
# bottom, left to right) so shared-string indices land in the same order
# as the authored workbook.
$ws.Range("F2").Value  = "value_added"
$ws.Range("G2").Value  = "integer"

$ws.Range("F3").Value  = "property_taxable_value"
$ws.Range("G3").Value  = "array"

$ws.Range("F4").Value  = "biodiesel_eq"
$ws.Range("G4").Value  = "integer"

$ws.Range("F5").Value  = "ethanol_eq"
$ws.Range("G5").Value  = "integer"

$ws.Range("F6").Value  = "fuel_taxable_value"
$ws.Range("G6").Value  = "array"

$ws.Range("F7").Value  = "property_taxable_value"
$ws.Range("G7").Value  = "array"

$ws.Range("F8").Value  = "NM_value"

$ws.Range("F9").Value  = "wages"

$ws.Range("F10").Value = "TCI"
$ws.Range("G10").Value = "integer"

$ws.Range("F11").Value = "TCI"
$ws.Range("G11").Value = "integer"

$ws.Range("F12").Value = "ethanol"
$ws.Range("G12").Value = "array"

$ws.Range("F13").Value = "TCI"
$ws.Range("G13").Value = "integer"

$ws.Range("F14").Value = "state_income_tax_assessed"
$ws.Range("G14").Value = "array"

$ws.Range("F15").Value = "ethanol"
$ws.Range("G15").Value = "array"

$ws.Range("F16").Value = "TCI"
$ws.Range("G16").Value = "integer"

$ws.Range("F17").Value = "TCI"
$ws.Range("G17").Value = "integer"

$ws.Range("F18").Value = "elec_eq"
$ws.Range("G18").Value = "integer"

$ws.Range("F19").Value = "state_income_tax_assessed"
$ws.Range("G19").Value = "array"

$ws.Range("F20").Value = "jobs_50"
$ws.Range("G20").Value = "integer"

$ws.Range("F21").Value = "ethanol"
$ws.Range("G21").Value = "array"

$ws.Range("F22").Value = "IA_value"
$ws.Range("G22").Value = "integer"

$ws.Range("F23").Value = "building_mats"
$ws.Range("G23").Value = "integer"

$ws.Range("F24").Value = "ethanol"
$ws.Range("G24").Value = "array"

# Row 25 only ever had a value in (what is now) column H - clear the blank
# placeholder cells Insert() left behind in F25:G25 so the row matches the
# original sparse layout.
$ws.Range("F25:G25").Clear()

# These two rows didn't get a "Variable type" - instead a note was added in
# column H (entered last, after every other Variable/Variable type cell
# above, which is why its shared-string id sorts after all the others).
$ws.Range("G8").Clear()
$ws.Range("G9").Clear()
$ws.Range("H8").Value = "determine whether to pass integer or array"
$ws.Range("H9").Value = "determine whether to pass integer or array"

# Approximate the author's resize of the two new columns to match column E's
# width.
$ws.Columns("F:G").ColumnWidth = 11.83

# Restore the cursor / selection to where the author left off.
$ws.Range("H15").Select()
